$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRows = @(5,6,10,11,12,13,14,17,20,21,22,23,25,26,27,29,30,31,32,33,35,37,38,39,40,41,43,45,47,48,49,51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '63.132.31'
$ws.Range("E2").Value = '  +5.67%  '
$ws.Range("D3").Value = '3.117.75'
$ws.Range("E3").Value = '  +3.63%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '584.53'
$ws.Range("E5").Value = '  +3.58%  '
$ws.Range("D6").Value = '144.58'
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.110.97'
$ws.Range("E8").Value = '  +3.84%  '
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  +11.57%  '
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  +8.44%  '
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  +6.88%  '
$ws.Range("D14").Value = '35.54'
$ws.Range("E14").Value = '  +4.64%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '3.637.48'
$ws.Range("E16").Value = '  +3.89%  '
$ws.Range("D17").Value = '7.18'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '3.118.17'
$ws.Range("E18").Value = '  +3.78%  '
$ws.Range("D19").Value = '63.072.87'
$ws.Range("E19").Value = '  +5.64%  '
$ws.Range("D20").Value = '466.63'
$ws.Range("E20").Value = '  +6.45%  '
$ws.Range("D21").Value = '14.05'
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("D22").Value = '0.726'
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("D23").Value = '7.52'
$ws.Range("E23").Value = '  +6.04%  '
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '82.00'
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '2.23'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("E28").Value = '  +4.99%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '8.25'
$ws.Range("E30").Value = '  +5.66%  '
$ws.Range("D31").Value = '6.84'
$ws.Range("E31").Value = '  +8.74%  '
$ws.Range("D32").Value = '26.99'
$ws.Range("E32").Value = '  +4.31%  '
$ws.Range("D33").Value = '0.110'
$ws.Range("E33").Value = '  +4.53%  '
$ws.Range("D34").Value = '0.0₃0867'
$ws.Range("E34").Value = '  +10.03%  '
$ws.Range("D35").Value = '2.38'
$ws.Range("E35").Value = '  +13.71%  '
$ws.Range("E36").Value = '  +4.96%  '
$ws.Range("D37").Value = '6.04'
$ws.Range("E37").Value = '  +2.19%  '
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  +16.16%  '
$ws.Range("D39").Value = '51.01'
$ws.Range("E39").Value = '  +3.91%  '
$ws.Range("D40").Value = '431.95'
$ws.Range("E40").Value = '  +7.32%  '
$ws.Range("D41").Value = '8.71'
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").Value = '2.928.10'
$ws.Range("E42").Value = '  +5.82%  '
$ws.Range("D43").Value = '0.0369'
$ws.Range("E43").Value = '  +4.51%  '
$ws.Range("E44").Value = '  +10.44%  '
$ws.Range("D45").Value = '0.112'
$ws.Range("E45").Value = '  +4.77%  '
$ws.Range("E46").Value = '  +6.04%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = '35.07'
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '123.64'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '24.50'
$ws.Range("E51").Value = '  +3.36%  '

foreach ($r in $textRows) {
    $ws.Range("D$r").Style = "Normal"
}
